$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing column B (Italian) entries whose text changed ---
$ws.Cells.Item(6, 2).Value = "Attestazione di identità "
$ws.Cells.Item(9, 2).Value = "Codice"

# --- Add new column C (English) and column D (German) translations ---
$ws.Cells.Item(2, 3).Value = "Request/application"
$ws.Cells.Item(2, 4).Value = "Gesuch / Anfrage"

$ws.Cells.Item(3, 3).Value = "Other documents"
$ws.Cells.Item(3, 4).Value = "Sonstige Dokumentation"

$ws.Cells.Item(4, 3).Value = "Payment declaration"
$ws.Cells.Item(4, 4).Value = "Zahlungsbestätigung"

$ws.Cells.Item(5, 3).Value = "Authorization Act"
$ws.Cells.Item(5, 4).Value = "Bewilligungsurkunde"

$ws.Cells.Item(6, 3).Value = "Identity Declaration"
$ws.Cells.Item(6, 4).Value = "Identitätsnachweis"

$ws.Cells.Item(7, 3).Value = "Administrative documentation"
$ws.Cells.Item(7, 4).Value = "Verwaltungsdokumentation"

$ws.Cells.Item(8, 3).Value = "Certification"
$ws.Cells.Item(8, 4).Value = "Bescheinigung"

$ws.Cells.Item(9, 3).Value = "Code"
$ws.Cells.Item(9, 4).Value = "Kode"

# --- Update header row (now 4 columns) ---
$ws.Cells.Item(1, 1).Value = "codice _1_livello"
$ws.Cells.Item(1, 2).Value = "label_ITA_1_livello"
$ws.Cells.Item(1, 3).Value = "label_ENG_1_livello"
$ws.Cells.Item(1, 4).Value = "label_DEU_1_livello"

# Header row is bold, like the original A1/B1 cells
$ws.Range("A1:D1").Font.Bold = $true

# Apply the distinct font used for the German (column D) translation cells
$ws.Range("D2:D9").Font.Name = "Calibri"

# --- Column widths for the two new columns ---
# (values chosen so the engine's internal pixel-rounding yields widths
# as close as possible to the target 26.1640625 / 24 character widths)
$ws.Range("C1").ColumnWidth = 25.3
$ws.Range("D1").ColumnWidth = 23.1

# --- Update selection / active cell ---
$ws.Range("D2").Select()
